$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: apply the "plain" text style (font Calibri 10 theme color 1,
# no fill, General number format, applyFont only) to a cell, matching the
# style used for most of the new data cells.
# ---------------------------------------------------------------------------
function Set-PlainStyle($cell) {
    $cell.Style = "Normal"
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

# ---------------------------------------------------------------------------
# Helper: apply the "plain+alignment" text style used for column F
# (same font, but also touches alignment so applyAlignment=1 is emitted).
# ---------------------------------------------------------------------------
function Set-PlainAlignStyle($cell) {
    $cell.Style = "Normal"
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
    $cell.WrapText = $false
}

# Scratch cell used to coerce numeric-looking text ("1994" etc.) into a
# real text/shared-string value instead of a number, without leaving any
# numeric display format behind on the real target cell.
$scratch = $ws.Range("AB1")
function Set-TextValue($cell, [string]$val) {
    $scratch.Formula = "'" + $val
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---------------------------------------------------------------------------
# New rows of data (MCH187 series)
# ---------------------------------------------------------------------------
$rows = @(
    @{ A = "MCH187-1"; C = "PUBLICATIONS, SOUTH AFRICA: APARTHEID IN CRISIS- A SPEECH BY RAYMOND SUTTNER, AUSTRALIAN ECONOMIC TIES WITH SA- MICHEAL LITTLE, AMANDLA (15 PUBLICATIONS), AFRICAN STUDIES REVIEW & NEWSLETTER- VOL. VVI NO.2, DEC 1994, VOLUME XVII NO. 1, JUNE 1995, WA CARE NEWS, SA ELECTIONS"; D = "1994"; DIsText = $true },
    @{ A = "MCH187-2"; C = "CARE NEWSLETTERS 1979, NO.11, 12, 13 1979, NO. 19, 26, 27 1980, NO.28 TO 36 1981"; D = "1979"; DIsText = $true },
    @{ A = "MCH187-3"; C = "CARE NEWSLETTERS NO. 37 TO 46 1982, NO. 47 TO 55 1983, NO. 56 TO 65 1984, NO. 66 TO 74 1984"; D = "1982"; DIsText = $true },
    @{ A = "MCH187-4"; C = "CARE NEWSLETTERS NO. 75 TO 80 1986, NO. 81 TO 85 1987, NO. 86 TO 94 1988, 1990"; D = "1986"; DIsText = $true },
    @{ A = "MCH187-5"; C = "VIVA NEWSLETTERS 1991-1993 "; D = "1991-1993"; DIsText = $false }
)

$E_VAL = "Series"
$F_VAL = "1 Box"
$G_VAL = "LOCATION: 23O | GRAP COUNT NUMER: NONE"

$r = 2
foreach ($row in $rows) {
    $cA = $ws.Cells.Item($r, 1)   # A
    $cC = $ws.Cells.Item($r, 3)   # C
    $cD = $ws.Cells.Item($r, 4)   # D
    $cE = $ws.Cells.Item($r, 5)   # E
    $cF = $ws.Cells.Item($r, 6)   # F
    $cG = $ws.Cells.Item($r, 7)   # G
    $cH = $ws.Cells.Item($r, 8)   # H

    Set-PlainStyle $cA
    $cA.Value2 = $row.A

    Set-PlainStyle $cC
    $cC.Value2 = $row.C

    Set-PlainStyle $cD
    if ($row.DIsText) {
        Set-TextValue $cD $row.D
    } else {
        $cD.Value2 = $row.D
    }

    Set-PlainStyle $cE
    $cE.Value2 = $E_VAL

    Set-PlainAlignStyle $cF
    $cF.Value2 = $F_VAL

    Set-PlainStyle $cG
    $cG.Value2 = $G_VAL

    Set-PlainStyle $cH

    $r++
}

$scratch.Clear() | Out-Null

# ---------------------------------------------------------------------------
# Restore the frozen-pane view and update the selection, as in the target.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:J6").Select() | Out-Null
